$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# Update the Runmode column for D Suite (row 5, now not running),
# E Suite (row 6, now running) and F Suite (row 7, now running).
$ws.Range("C5").Value = "N"
$ws.Range("C6").Value = "Y"
$ws.Range("C7").Value = "Y"

# Move the active selection to C7 to match the saved view state.
$ws.Range("C7").Select()
